# Apply updated coin price / volume(1h) values from the Feb 10 2023 symbol-list refresh.
# Cells D (Price) and E (Volume(1h)) hold plain text (e.g. "307.19", "-4.57%"), not
# numbers, in the source sheet. We force the Text number format before assigning the
# new value so Excel stores it as a string rather than auto-coercing it to a numeric
# value, then clear formatting again so the cell keeps its original (default) style.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rng = $ws.Range("D2:E2")
$rng.NumberFormat = "@"
$ws.Range("D2").Value = "307.19"
$ws.Range("E2").Value = "-4.57%"
$rng.ClearFormats()

$rng = $ws.Range("D3:E3")
$rng.NumberFormat = "@"
$ws.Range("D3").Value = "39.14"
$ws.Range("E3").Value = "-8.84%"
$rng.ClearFormats()

$rng = $ws.Range("D4:E4")
$rng.NumberFormat = "@"
$ws.Range("D4").Value = "5.100"
$ws.Range("E4").Value = "-1.95%"
$rng.ClearFormats()

$rng = $ws.Range("D5:E5")
$rng.NumberFormat = "@"
$ws.Range("D5").Value = "0.07679"
$ws.Range("E5").Value = "-6.09%"
$rng.ClearFormats()

$rng = $ws.Range("D6:E6")
$rng.NumberFormat = "@"
$ws.Range("D6").Value = "4.238"
$ws.Range("E6").Value = "-1.81%"
$rng.ClearFormats()

$rng = $ws.Range("D7:E7")
$rng.NumberFormat = "@"
$ws.Range("D7").Value = "1.603"
$ws.Range("E7").Value = "-11.17%"
$rng.ClearFormats()

$rng = $ws.Range("D8:E8")
$rng.NumberFormat = "@"
$ws.Range("D8").Value = "0.9162"
$ws.Range("E8").Value = "-3.48%"
$rng.ClearFormats()

$rng = $ws.Range("E9")
$rng.NumberFormat = "@"
$ws.Range("E9").Value = "-7.70%"
$rng.ClearFormats()

$rng = $ws.Range("D10:E10")
$rng.NumberFormat = "@"
$ws.Range("D10").Value = "0.1751"
$ws.Range("E10").Value = "-6.81%"
$rng.ClearFormats()

$rng = $ws.Range("D11:E11")
$rng.NumberFormat = "@"
$ws.Range("D11").Value = "0.09008"
$ws.Range("E11").Value = "-3.79%"
$rng.ClearFormats()

$rng = $ws.Range("D12:E12")
$rng.NumberFormat = "@"
$ws.Range("D12").Value = "0.04427"
$ws.Range("E12").Value = "-4.15%"
$rng.ClearFormats()

$rng = $ws.Range("E13")
$rng.NumberFormat = "@"
$ws.Range("E13").Value = "-0.33%"
$rng.ClearFormats()

$rng = $ws.Range("D14:E14")
$rng.NumberFormat = "@"
$ws.Range("D14").Value = "0.001259"
$ws.Range("E14").Value = "-2.23%"
$rng.ClearFormats()

$rng = $ws.Range("D15:E15")
$rng.NumberFormat = "@"
$ws.Range("D15").Value = "0.005812"
$ws.Range("E15").Value = "0.87%"
$rng.ClearFormats()

$rng = $ws.Range("E16")
$rng.NumberFormat = "@"
$ws.Range("E16").Value = "2,416.50%"
$rng.ClearFormats()

$rng = $ws.Range("E17")
$rng.NumberFormat = "@"
$ws.Range("E17").Value = "-0.08%"
$rng.ClearFormats()

$rng = $ws.Range("E18")
$rng.NumberFormat = "@"
$ws.Range("E18").Value = "-4.69%"
$rng.ClearFormats()

$rng = $ws.Range("D19:E19")
$rng.NumberFormat = "@"
$ws.Range("D19").Value = "0.3305"
$ws.Range("E19").Value = "-1.73%"
$rng.ClearFormats()

$rng = $ws.Range("D20:E20")
$rng.NumberFormat = "@"
$ws.Range("D20").Value = "7.039"
$ws.Range("E20").Value = "-5.64%"
$rng.ClearFormats()

$rng = $ws.Range("D21:E21")
$rng.NumberFormat = "@"
$ws.Range("D21").Value = "0.1346"
$ws.Range("E21").Value = "-3.03%"
$rng.ClearFormats()

$rng = $ws.Range("D22:E22")
$rng.NumberFormat = "@"
$ws.Range("D22").Value = "0.2758"
$ws.Range("E22").Value = "8.27%"
$rng.ClearFormats()

$rng = $ws.Range("D23:E23")
$rng.NumberFormat = "@"
$ws.Range("D23").Value = "0.04138"
$ws.Range("E23").Value = "-0.23%"
$rng.ClearFormats()

$rng = $ws.Range("E24")
$rng.NumberFormat = "@"
$ws.Range("E24").Value = "-3.56%"
$rng.ClearFormats()

$rng = $ws.Range("D25:E25")
$rng.NumberFormat = "@"
$ws.Range("D25").Value = "0.004110"
$ws.Range("E25").Value = "-4.11%"
$rng.ClearFormats()

$rng = $ws.Range("D26:E26")
$rng.NumberFormat = "@"
$ws.Range("D26").Value = "0.0001301"
$ws.Range("E26").Value = "8.51%"
$rng.ClearFormats()

$rng = $ws.Range("E38")
$rng.NumberFormat = "@"
$ws.Range("E38").Value = "-10.41%"
$rng.ClearFormats()

$rng = $ws.Range("D39:E39")
$rng.NumberFormat = "@"
$ws.Range("D39").Value = "0.05200"
$ws.Range("E39").Value = "-6.60%"
$rng.ClearFormats()

$rng = $ws.Range("D40:E40")
$rng.NumberFormat = "@"
$ws.Range("D40").Value = "0.007922"
$ws.Range("E40").Value = "-2.27%"
$rng.ClearFormats()

$rng = $ws.Range("D41:E41")
$rng.NumberFormat = "@"
$ws.Range("D41").Value = "0.1318"
$ws.Range("E41").Value = "-5.89%"
$rng.ClearFormats()

$rng = $ws.Range("D42:E42")
$rng.NumberFormat = "@"
$ws.Range("D42").Value = "0.007206"
$ws.Range("E42").Value = "10.14%"
$rng.ClearFormats()

$rng = $ws.Range("D43:E43")
$rng.NumberFormat = "@"
$ws.Range("D43").Value = "0.001951"
$ws.Range("E43").Value = "-7.49%"
$rng.ClearFormats()

$rng = $ws.Range("D44:E44")
$rng.NumberFormat = "@"
$ws.Range("D44").Value = "0.008393"
$ws.Range("E44").Value = "9.65%"
$rng.ClearFormats()

$rng = $ws.Range("D45")
$rng.NumberFormat = "@"
$ws.Range("D45").Value = "0.3335"
$rng.ClearFormats()

$rng = $ws.Range("D46:E46")
$rng.NumberFormat = "@"
$ws.Range("D46").Value = "0.00006427"
$ws.Range("E46").Value = "-4.58%"
$rng.ClearFormats()

$rng = $ws.Range("E47")
$rng.NumberFormat = "@"
$ws.Range("E47").Value = "0.16%"
$rng.ClearFormats()

$rng = $ws.Range("E48")
$rng.NumberFormat = "@"
$ws.Range("E48").Value = "-26.71%"
$rng.ClearFormats()

$rng = $ws.Range("D49:E49")
$rng.NumberFormat = "@"
$ws.Range("D49").Value = "0.004241"
$ws.Range("E49").Value = "38.06%"
$rng.ClearFormats()

$rng = $ws.Range("D50:E50")
$rng.NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").Value = "0.16%"
$rng.ClearFormats()

$rng = $ws.Range("E51")
$rng.NumberFormat = "@"
$ws.Range("E51").Value = "0.16%"
$rng.ClearFormats()
